$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The source data already has one row per year (2018, 2019, 2020 -> rows 2-4).
# Append the 2021 figures as the next row (row 5), reusing the same column
# layout as the existing years.
$ws.Range("A5").Value = '2021年'
$ws.Range("B5").Value = 8.4
$ws.Range("C5").Value = 25.5
$ws.Range("D5").Value = -26.6
$ws.Range("E5").Value = -21.1
$ws.Range("F5").Value = -4.3
$ws.Range("G5").Formula = '=""'
$ws.Range("H5").Value = 13
$ws.Range("I5").Value = -5.4
$ws.Range("J5").Value = -6.8
$ws.Range("K5").Value = 13.1
$ws.Range("L5").Value = -0.8
$ws.Range("M5").Value = 40.5
$ws.Range("N5").Value = -25.4
$ws.Range("O5").Value = -3.2
$ws.Range("P5").Value = 29.6
$ws.Range("Q5").Value = -12.8
$ws.Range("R5").Formula = '=""'
$ws.Range("S5").Value = 159.9
$ws.Range("T5").Value = 3.8
$ws.Range("U5").Value = 2.3
$ws.Range("V5").Value = -1.3
$ws.Range("W5").Value = -7
$ws.Range("X5").Value = 11.2
$ws.Range("Y5").Value = 29.2
$ws.Range("Z5").Value = 221.5
$ws.Range("AA5").Value = 35.5
$ws.Range("AB5").Value = 5.9
$ws.Range("AC5").Value = 8.9
$ws.Range("AD5").Value = 34.2
$ws.Range("AE5").Value = -13.2
$ws.Range("AF5").Value = -25.6
$ws.Range("AG5").Formula = '=""'
$ws.Range("AH5").Value = -29.8
$ws.Range("AI5").Value = 16.7
$ws.Range("AJ5").Value = -52.8
$ws.Range("AK5").Value = 55.1
$ws.Range("AL5").Value = 6
$ws.Range("AM5").Value = -24.5
$ws.Range("AN5").Value = -15.8
$ws.Range("AO5").Value = -23.3
$ws.Range("AP5").Value = 36.8
$ws.Range("AQ5").Value = 2.5
$ws.Range("AR5").Value = -10.3
$ws.Range("AS5").Formula = '=""'
$ws.Range("AT5").Value = 231.7
$ws.Range("AU5").Formula = '=""'
$ws.Range("AV5").Value = 9.1
$ws.Range("AW5").Value = -44.4
$ws.Range("AX5").Value = -19.6
$ws.Range("AY5").Value = -15.8
$ws.Range("AZ5").Formula = '=""'
$ws.Range("BA5").Value = 0.5
$ws.Range("BB5").Value = -15.5
$ws.Range("BC5").Value = 29.4
$ws.Range("BD5").Formula = '=""'
$ws.Range("BE5").Value = 16.1
$ws.Range("BF5").Value = 102.9
$ws.Range("BG5").Value = 17.7
$ws.Range("BH5").Value = 57
$ws.Range("BI5").Value = 1.8
$ws.Range("BJ5").Value = 38.9
$ws.Range("BK5").Value = 29.1
$ws.Range("BL5").Value = -5.6
$ws.Range("BM5").Value = -9.1
$ws.Range("BN5").Value = -11.6
$ws.Range("BO5").Value = -33.3
$ws.Range("BP5").Value = 81.5
$ws.Range("BQ5").Value = 103
$ws.Range("BR5").Value = -6.5
$ws.Range("BS5").Value = -23.3
$ws.Range("BT5").Value = -17.3
$ws.Range("BU5").Value = -26.4
$ws.Range("BV5").Value = -5.7
$ws.Range("BW5").Value = -0.4
$ws.Range("BX5").Value = -5
$ws.Range("BY5").Value = -0.3
$ws.Range("BZ5").Value = -10.4
$ws.Range("CA5").Value = -57.4
$ws.Range("CB5").Value = 3.9
$ws.Range("CC5").Value = -29.7
$ws.Range("CD5").Value = -3.7
$ws.Range("CE5").Value = -15.1
$ws.Range("CF5").Value = -4.2
$ws.Range("CG5").Value = 15.5
$ws.Range("CH5").Value = 39.4
$ws.Range("CI5").Value = -10.8
$ws.Range("CJ5").Value = -17.4
$ws.Range("CK5").Value = 10.9
$ws.Range("CL5").Value = 48.6
$ws.Range("CM5").Value = 90.8
$ws.Range("CN5").Value = 0.8
$ws.Range("CO5").Value = 19.8
$ws.Range("CP5").Value = 29.1
$ws.Range("CQ5").Value = 64.3
$ws.Range("CR5").Value = 83
$ws.Range("CS5").Value = 197.2
$ws.Range("CT5").Value = 15.8
$ws.Range("CU5").Value = 22.2
$ws.Range("CV5").Value = 1
$ws.Range("CW5").Value = -70
$ws.Range("CX5").Value = 27.8
$ws.Range("CY5").Value = 11.1
$ws.Range("CZ5").Value = -21
$ws.Range("DA5").Value = 3.4
$ws.Range("DB5").Value = 70.9
$ws.Range("DC5").Value = 0.7
$ws.Range("DD5").Value = -31.7
$ws.Range("DE5").Value = -13.8
$ws.Range("DF5").Value = 22.2
$ws.Range("DG5").Value = 58.7
$ws.Range("DH5").Value = 0.2
$ws.Range("DI5").Value = -20.6
$ws.Range("DJ5").Value = 1.4
$ws.Range("DK5").Value = -69.3
# Match the formatting used for the year label in column A of the previous
# rows (bold, bordered, centered header-style cell).
$ws.Range("A4").Copy() | Out-Null
$ws.Range("A5").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
